$d = $word.ActiveDocument

# --- First paragraph: "**ID__AFFARS_5308_topic_4__ID** " ---
$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right, 5-twip space, no line) and
# bump the left indent from 120 -> 225 twips (225 twips = 11.25 pt).
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5
$p1.Format.LeftIndent = 11.25

# Drop the trailing " " run (second run in the paragraph, right before the
# paragraph mark) while leaving the first run's formatting untouched.
$spaceRange = $d.Range($p1.Range.End - 2, $p1.Range.End - 1)
$spaceRange.Delete()

# Update the bookmark-style placeholder text carried by the (now sole) run.
$d.Content.Find.Execute("**ID__AFFARS_5308_topic_4__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5308_405_3__ID**", 2)
